$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (columns A:R, data ends around row 308, totals row 309) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new data row at position 165 (shifts everything from old row 165
# down to row 166, ..., old totals row 308 becomes row 309).
$ws1.Rows("165:165").Insert()

# Populate the newly inserted row with the new client entry.
$ws1.Cells.Item(165, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Cells.Item(165, 2).Value = "BRITO CARDENAS RUTH CECILIA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(165, $c).Value = 0
}

# Fix up the "x de 306" -> "x de 307" counters on the (now shifted) totals row 309.
$ws1.Cells.Item(309, 3).Value  = "0 de 307"
$ws1.Cells.Item(309, 4).Value  = "1 de 307"
$ws1.Cells.Item(309, 5).Value  = "1 de 307"
$ws1.Cells.Item(309, 6).Value  = "0 de 307"
$ws1.Cells.Item(309, 7).Value  = "0 de 307"
$ws1.Cells.Item(309, 8).Value  = "2 de 307"
$ws1.Cells.Item(309, 9).Value  = "1 de 307"
$ws1.Cells.Item(309, 10).Value = "0 de 307"
$ws1.Cells.Item(309, 11).Value = "0 de 307"
$ws1.Cells.Item(309, 12).Value = "1 de 307"
$ws1.Cells.Item(309, 13).Value = "6 de 307"
$ws1.Cells.Item(309, 14).Value = "0 de 307"
$ws1.Cells.Item(309, 15).Value = "0 de 307"
$ws1.Cells.Item(309, 16).Value = "0 de 307"
$ws1.Cells.Item(309, 17).Value = "0 de 307"
$ws1.Cells.Item(309, 18).Value = "0 de 307"

# --- Sheet "VENTA MENSUAL" (columns A:G, data ends around row 308, totals row 309) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same single-row insert of the new client entry.
$ws2.Rows("165:165").Insert()

$ws2.Cells.Item(165, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Cells.Item(165, 2).Value = "BRITO CARDENAS RUTH CECILIA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(165, $c).Value = 0
}
# The totals row on this sheet holds plain numeric sums (unaffected by the
# new all-zero row), so it needs no further edits after the shift.
